# Auto update stock data
# Updates the report date from 2025/12/27 to 2025/12/28 in column A
# (Date_1) for every stock block on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$targetRows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $targetRows) {
    $cell = $ws.Cells.Item($r, 1)
    # Force the cell to stay a text value (it's stored as text, not a
    # real date) instead of letting Excel auto-convert the
    # yyyy/mm/dd-looking string into a date serial number.
    $cell.NumberFormat = "@"
    $cell.Value = "2025/12/28"
    # Restore the default ("Normal") style so no stray number format is
    # left attached to the cell, matching the original formatting.
    $cell.Style = "Normal"
}
